$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BitácoraExperiencia1")
$ws.Activate()

# Update the "Equipo" number from 1 to 9
$ws.Range("C2").Value = 9

# Delete the last activity row (row 18): 12 | envie el proyecto a la profesora | Aaron Tobar
$ws.Rows.Item(18).Delete()

# Scroll so row 9 is at the top of the visible window, then select B18,
# matching the author's final on-screen view state.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B18").Select()
